$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "https://raw.githubusercontent.com/bryanparthum/farmland_conservation/master/images/"
$newPrefix = "https://bryanparthum.github.io/farmland_conservation/choice_cards/cards/"

# Card URL columns are Q..W (17..23); data rows are 2..85.
for ($r = 2; $r -le 85; $r++) {
    for ($c = 17; $c -le 23; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = [string]$cell.Value2
        if ($val -ne $null -and $val.StartsWith($oldPrefix)) {
            $cell.Value = $val.Replace($oldPrefix, $newPrefix)
        }
    }
}
